$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "I am right on schedule to finish the " ->
#           "I am right on schedule to finish "
#
# This text lives in a run that is immediately followed by another run
# ("to finish the project and present for next week. ...") which, by
# coincidence, carries exactly the same run formatting. A plain
# Find/Replace (or a Range.Text / Range.Delete on the "the " text)
# causes the two neighbouring, identically-formatted runs to be
# re-merged into a single <w:r>, which is not what happened in the real
# edit (the two runs stay distinct there). Temporarily wrapping just the
# doomed "the " in a bookmark and clearing the bookmark's own Range
# edits that run in isolation, leaving its neighbour run untouched.
# ---------------------------------------------------------------------
$anchor = "schedule to finish "
$deadWord = "the "

$r = $d.Content
$r.Find.Execute($anchor + $deadWord, $true, $false, $false, $false, $false, `
                $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.MoveStart(1, -1 * $deadWord.Length) | Out-Null

$d.Bookmarks.Add("TmpDeleteMarker", $r) | Out-Null
$d.Bookmarks("TmpDeleteMarker").Range.Text = ""
$d.Bookmarks("TmpDeleteMarker").Delete()

# ---------------------------------------------------------------------
# Change 2: the hidden "_GoBack" bookmark (Word's "last edit position"
# marker) moves from right after "...production testing." (end of the
# third body paragraph) to a collapsed position in the middle of the
# word "focussed" ("...Therefore, I fo|cussed solely...") inside the
# "Financially, ..." paragraph.
#
# Word only ever keeps a single "_GoBack" bookmark in a document, so
# re-adding a bookmark under that same name both deletes the old one
# and inserts the new one, in the right place, in one step.
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("I fo", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2) | Out-Null

Write-Output "edits applied"
